$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.548.43"
$ws.Range("E2").Value = "  -2.74%  "
$ws.Range("D3").Value = "2.975.53"
$ws.Range("E3").Value = "  -4.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "497.00"
$ws.Range("E5").Value = "  -4.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.56"
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "2.973.62"
$ws.Range("E8").Value = "  -4.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.425"
$ws.Range("E9").Value = "  -3.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.27"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("E12").Value = "  -7.68%  "
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "3.490.22"
$ws.Range("E14").Value = "  -4.71%  "
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("D16").Value = "56.346.27"
$ws.Range("E16").Value = "  -2.85%  "
$ws.Range("D17").Value = "2.979.99"
$ws.Range("E17").Value = "  -4.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000146"
$ws.Range("E18").Value = "  -3.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.81"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.23"
$ws.Range("E20").Value = "  -5.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.71"
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "322.94"
$ws.Range("E22").Value = "  -6.03%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  -8.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "61.41"
$ws.Range("E25").Value = "  -10.62%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  -1.42%  "
$ws.Range("D28").Value = "0.0₃0897"
$ws.Range("E28").Value = "  -5.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.52"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.77"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.73"
$ws.Range("E33").Value = "  -6.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.80"
$ws.Range("E34").Value = "  -8.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.47"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.46"
$ws.Range("E36").Value = "  -6.22%  "
$ws.Range("E37").Value = "  -5.65%  "
$ws.Range("E38").Value = "  -8.98%  "
$ws.Range("E39").Value = "  -3.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.24"
$ws.Range("E40").Value = "  -3.68%  "
$ws.Range("D41").Value = "3.005.59"
$ws.Range("E41").Value = "  -4.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.24"
$ws.Range("E42").Value = "  -7.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  -6.74%  "
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.634"
$ws.Range("E46").Value = "  -8.81%  "
$ws.Range("D47").Value = "2.207.19"
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("E48").Value = "  -8.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.96"
$ws.Range("E49").Value = "  +8.46%  "
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.16"
$ws.Range("E51").Value = "  -5.94%  "
